# Apply cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.417.88"
$ws.Range("E2").Value = "  -5.14%  "

# Row 3
$ws.Range("D3").Value = "3.466.58"
$ws.Range("E3").Value = "  -6.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.23"
$ws.Range("E5").Value = "  -7.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.92"
$ws.Range("E6").Value = "  -8.81%  "

# Row 7
$ws.Range("D7").Value = "3.469.54"
$ws.Range("E7").Value = "  -6.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  -5.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  -7.23%  "

# Row 11
$ws.Range("E11").Value = "  -5.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -6.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000216"
$ws.Range("E13").Value = "  -8.25%  "

# Row 14
$ws.Range("D14").Value = "4.059.62"
$ws.Range("E14").Value = "  -5.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.08"
$ws.Range("E15").Value = "  -5.83%  "

# Row 16
$ws.Range("D16").Value = "3.507.00"
$ws.Range("E16").Value = "  -5.16%  "

# Row 17
$ws.Range("D17").Value = "66.479.77"
$ws.Range("E17").Value = "  -4.92%  "

# Row 18
$ws.Range("E18").Value = "  -0.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  -4.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").Value = "  -7.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "438.93"
$ws.Range("E21").Value = "  -7.56%  "

# Row 22
$ws.Range("E22").Value = "  -15.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.617"
$ws.Range("E23").Value = "  -5.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.68"
$ws.Range("E24").Value = "  -4.26%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("D26").Value = "3.619.25"
$ws.Range("E26").Value = "  -5.70%  "

# Row 27
$ws.Range("E27").Value = "  -6.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -9.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  -12.91%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.49"
$ws.Range("E30").Value = "  -6.70%  "

# Row 31
$ws.Range("E31").Value = "  -10.41%  "

# Row 32
$ws.Range("E32").Value = "  +0.10%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.26"
$ws.Range("E33").Value = "  -6.49%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  -5.92%  "

# Row 35
$ws.Range("E35").Value = "  -8.12%  "

# Row 36
$ws.Range("D36").Value = "3.463.70"
$ws.Range("E36").Value = "  -6.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -9.99%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.86"
$ws.Range("E38").Value = "  -7.48%  "

# Row 39
$ws.Range("E39").Value = "  +0.14%  "

# Row 40
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.36"
$ws.Range("E41").Value = "  -4.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.14"
$ws.Range("E42").Value = "  -6.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  -8.62%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0849"
$ws.Range("E44").Value = "  -6.67%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.876"
$ws.Range("E45").Value = "  -6.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.98"
$ws.Range("E46").Value = "  -4.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.87"
$ws.Range("E47").Value = "  -8.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  -5.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.45"
$ws.Range("E49").Value = "  -5.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.42"
$ws.Range("E50").Value = "  -13.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.990"
$ws.Range("E51").Value = "  -7.80%  "
